$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format so numeric-looking price strings (e.g. "1.00")
# are preserved as text instead of being coerced to numbers, matching the source data.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "46.199.58"
$ws.Range("E2").Value = "  -1.58%  "

# Row 3
$ws.Range("D3").Value = "2.356.61"
$ws.Range("E3").Value = "  +1.86%  "

# Row 4
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.10%  "

# Row 5
$ws.Range("D5").Value = "301.76"
$ws.Range("E5").Value = "  +0.97%  "

# Row 6
$ws.Range("D6").Value = "99.50"
$ws.Range("E6").Value = "  +0.51%  "

# Row 7
$ws.Range("D7").Value = "0.570"
$ws.Range("E7").Value = "  -0.42%  "

# Row 8
$ws.Range("E8").Value = "  +0.11%  "

# Row 9
$ws.Range("E9").Value = "  -3.19%  "

# Row 10
$ws.Range("D10").Value = "34.43"
$ws.Range("E10").Value = "  -3.94%  "

# Row 11
$ws.Range("E11").Value = "  -0.06%  "

# Row 12
$ws.Range("D12").Value = "7.13"
$ws.Range("E12").Value = "  -2.91%  "

# Row 13
$ws.Range("D13").Value = "0.103"
$ws.Range("E13").Value = "  -0.39%  "

# Row 14
$ws.Range("D14").Value = "2.715.83"
$ws.Range("E14").Value = "  +1.97%  "

# Row 15
$ws.Range("D15").Value = "2.353.29"
$ws.Range("E15").Value = "  +1.92%  "

# Row 16
$ws.Range("E16").Value = "  -0.85%  "

# Row 17
$ws.Range("D17").Value = "13.62"
$ws.Range("E17").Value = "  -2.66%  "

# Row 18
$ws.Range("D18").Value = "46.112.65"
$ws.Range("E18").Value = "  -1.41%  "

# Row 19
$ws.Range("D19").Value = "12.72"
$ws.Range("E19").Value = "  -3.10%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0964"
$ws.Range("E20").Value = "  +2.52%  "

# Row 21
$ws.Range("E21").Value = "  -1.13%  "

# Row 23
$ws.Range("D23").Value = "246.64"
$ws.Range("E23").Value = "  -0.93%  "

# Row 24
$ws.Range("D24").Value = "2.83"
$ws.Range("E24").Value = "  -3.19%  "

# Row 25
$ws.Range("E25").Value = "  -0.18%  "

# Row 26
$ws.Range("E26").Value = "  -3.26%  "

# Row 27
$ws.Range("D27").Value = "39.54"
$ws.Range("E27").Value = "  -7.71%  "

# Row 28
$ws.Range("D28").Value = "2.20"
$ws.Range("E28").Value = "  -2.85%  "

# Row 29
$ws.Range("E29").Value = "  -0.68%  "

# Row 30
$ws.Range("E30").Value = "  +4.03%  "

# Row 31
$ws.Range("D31").Value = "3.73"
$ws.Range("E31").Value = "  +19.67%  "

# Row 32
$ws.Range("E32").Value = "  +5.86%  "

# Row 33
$ws.Range("D33").Value = "5.53"
$ws.Range("E33").Value = "  -3.86%  "

# Row 34
$ws.Range("D34").Value = "145.95"
$ws.Range("E34").Value = "  -1.03%  "

# Row 35
$ws.Range("E35").Value = "  -3.14%  "

# Row 36
$ws.Range("E36").Value = "  -1.32%  "

# Row 37
$ws.Range("D37").Value = "1.89"
$ws.Range("E37").Value = "  +4.93%  "

# Row 38
$ws.Range("E38").Value = "  -2.07%  "

# Row 39
$ws.Range("D39").Value = "14.94"
$ws.Range("E39").Value = "  -5.40%  "

# Row 40
$ws.Range("E40").Value = "  -1.17%  "

# Row 41
$ws.Range("E41").Value = "  -2.33%  "

# Row 42
$ws.Range("D42").Value = "3.22"
$ws.Range("E42").Value = "  -6.17%  "

# Row 43
$ws.Range("D43").Value = "1.889.65"
$ws.Range("E43").Value = "  +2.62%  "

# Row 44
$ws.Range("D44").Value = "1.00"
$ws.Range("E44").Value = "  +0.03%  "

# Row 45
$ws.Range("D45").Value = "93.11"
$ws.Range("E45").Value = "  +2.22%  "

# Row 46
$ws.Range("E46").Value = "  -9.56%  "

# Row 47
$ws.Range("E47").Value = "  -6.31%  "

# Row 48
$ws.Range("D48").Value = "8.24"
$ws.Range("E48").Value = "  +3.23%  "

# Row 49
$ws.Range("D49").Value = "97.71"
$ws.Range("E49").Value = "  +0.38%  "

# Row 50
$ws.Range("D50").Value = "2.586.88"
$ws.Range("E50").Value = "  +1.78%  "

# Row 51: coin replaced (ordi -> EnergySwap)
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "14.50"
$ws.Range("E51").Value = "  +5.73%  "

# Remove the temporary text number-format so the cell style matches the original
# (unstyled) cells rather than leaving a stray "@" format behind.
$ws.Range("D2:D51").ClearFormats()
